# Discharge_Aug02.xlsx -- "lots of discharge data"
#
# Adds a third depth/velocity/discharge sub-table (rows 32-44) to the
# "stn3" sheet, mirroring the sub-table that already exists on "stn1"
# (rows 29-40), and leaves the workbook with "stn3" as the active sheet.

$wb = $excel.ActiveWorkbook

$stn1 = $wb.Worksheets.Item("stn1")
$stn3 = $wb.Worksheets.Item("stn3")

# ---------------------------------------------------------------------
# stn3 ("xl/worksheets/sheet2.xml"): new header + data rows 32-44
# ---------------------------------------------------------------------

# Row 32: bold "New depth" label (black font color, like the sample in the
# source file -- a fresh font/cellXf is created for this).
$stn3.Range("A32").Value = "New depth"
$stn3.Range("A32").Font.Bold = $true
$stn3.Range("A32").Font.Color = 0

# Row 33: column headers (reuses existing shared strings X / V / D /
# segment / Q / Qtotal).
$stn3.Range("A33").Value = "X"
$stn3.Range("B33").Value = "V"
$stn3.Range("C33").Value = "D"
$stn3.Range("D33").Value = "segment"
$stn3.Range("E33").Value = "Q"
$stn3.Range("F33").Value = "Qtotal"

# Row 34: first data row -- literal values copied from A18/B18, a fresh
# C34 formula (cm conversion of C18) and the running D/F formulas.
$stn3.Range("A34").Value = 0.5
$stn3.Range("B34").Value = 0
$stn3.Range("C34").Formula = "=C18*2.54"
$stn3.Range("D34").Formula = "=A34"
$stn3.Range("F34").Formula = "=SUM(E34:E44)"

# Rows 35-44: A/B are literal values copied from A19:A28 / B19:B28, C is
# the cm conversion of C19:C28, D is the running midpoint formula and E is
# the segment-discharge formula.
$aVals = @(0.55, 0.6, 0.65, 0.7, 0.75, 0.8, 0.85, 0.9, 0.95, 1)
$bVals = @(0.13156, 0.18304, 0.12584, 0.0572, 0.0572, 0.04576, 0.04004, 0.04004, 0, 0)

for ($i = 0; $i -lt 10; $i++) {
    $row = 35 + $i
    $srcRow = 19 + $i

    $stn3.Range("A$row").Value = $aVals[$i]
    $stn3.Range("B$row").Value = $bVals[$i]
    $stn3.Range("C$row").Formula = "=C$srcRow*2.54"

    if ($row -eq 44) {
        $stn3.Range("D$row").Formula = "=A$row"
    } else {
        $nextRow = $row + 1
        $stn3.Range("D$row").Formula = "=(A$row+(A$nextRow-A$row)/2)"
    }

    $prevRow = $row - 1
    $stn3.Range("E$row").Formula = "=(D$row-D$prevRow)*(B$row)*C$row"
}

# ---------------------------------------------------------------------
# Selections + active sheet
# ---------------------------------------------------------------------

$stn1.Activate()
$stn1.Range("C33").Select()

$stn3.Activate()
$stn3.Range("F34").Select()
